# Adds three new daily rows (06, 07 and 08-10-2021) to the bottom of the table,
# matching the MV -datos- source update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Set-DateText($row, $text) {
    # The column-A values look like dates (dd-mm-yyyy) and Excel would normally
    # auto-convert a literal string assignment into a date serial number. Enter
    # it as a text formula instead and then collapse it down to a plain value via
    # copy / paste-special so the cell ends up as a regular shared-string cell
    # (matching the un-styled string cells used for every other date in column A).
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = "=""$text"""
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

# ---- Row 193 (06-10-2021) ----
Set-DateText 193 "06-10-2021"
$ws.Cells.Item(193, 2).Value = 23.97
$ws.Cells.Item(193, 3).Value = 810.63
$ws.Cells.Item(193, 4).Value = 194.62
$ws.Cells.Item(193, 5).Value = 118.17
$ws.Cells.Item(193, 6).Value = 1.29
$ws.Cells.Item(193, 7).Value = 37.12
$ws.Cells.Item(193, 8).Value = 126.46
$ws.Cells.Item(193, 9).Value = 6.33
$ws.Cells.Item(193, 10).Value = 95.06
$ws.Cells.Item(193, 11).Value = 92.85
$ws.Cells.Item(193, 12).Value = 1145.12
$ws.Cells.Item(193, 13).Value = 220.76
$ws.Cells.Item(193, 14).Value = 89.67
$ws.Cells.Item(193, 15).Value = 591.44
$ws.Cells.Item(193, 16).Value = 644.84
$ws.Cells.Item(193, 17).Value = 810.63
$ws.Cells.Item(193, 18).Value = 810.63
$ws.Cells.Item(193, 19).Value = 385.46
$ws.Cells.Item(193, 20).Value = 988.57
$ws.Cells.Item(193, 21).Value = 810.63
$ws.Cells.Item(193, 22).Value = 597.55
$ws.Cells.Item(193, 23).Value = 104.13
$ws.Cells.Item(193, 24).Value = 564.94
$ws.Cells.Item(193, 25).Value = 29.07
$ws.Cells.Item(193, 26).Value = 0.04
$ws.Cells.Item(193, 27).Value = 940.62
$ws.Cells.Item(193, 28).Value = 2.63
$ws.Cells.Item(193, 29).Value = 7.92
$ws.Cells.Item(193, 30).Value = 873.81
$ws.Cells.Item(193, 31).Value = 0.12
$ws.Cells.Item(193, 32).Value = 30.71
$ws.Cells.Item(193, 33).Value = 190.29
$ws.Cells.Item(193, 34).Value = 51.73
$ws.Cells.Item(193, 35).Value = 1105.3
$ws.Cells.Item(193, 36).Value = 91.38
$ws.Cells.Item(193, 37).Value = 196.35
$ws.Cells.Item(193, 38).Value = 8.2
$ws.Cells.Item(193, 39).Value = 0.21
$ws.Cells.Item(193, 40).Value = 33.78
$ws.Cells.Item(193, 41).Value = 14.43
$ws.Cells.Item(193, 42).Value = 16.01
$ws.Cells.Item(193, 43).Value = 39.48
$ws.Cells.Item(193, 44).Value = 18.92
$ws.Cells.Item(193, 45).Value = 104.91
$ws.Cells.Item(193, 46).Value = 54.09
$ws.Cells.Item(193, 47).Value = 148.25
$ws.Cells.Item(193, 48).Value = 0.02
$ws.Cells.Item(193, 49).Value = 216.17
$ws.Cells.Item(193, 50).Value = 194
$ws.Cells.Item(193, 51).Value = 11.21
$ws.Cells.Item(193, 52).Value = 10.89
$ws.Cells.Item(193, 53).Value = 0.06
$ws.Cells.Item(193, 54).Value = 4.75
$ws.Cells.Item(193, 55).Value = 250.92
$ws.Cells.Item(193, 56).Value = 1.91
$ws.Cells.Item(193, 57).Value = 0.68
$ws.Cells.Item(193, 58).Value = 7.27
$ws.Cells.Item(193, 59).Value = 125.77
$ws.Cells.Item(193, 60).Value = 204.02

# ---- Row 194 (07-10-2021) ----
Set-DateText 194 "07-10-2021"
$ws.Cells.Item(194, 2).Value = 24.14
$ws.Cells.Item(194, 3).Value = 816.28
$ws.Cells.Item(194, 4).Value = 195.1
$ws.Cells.Item(194, 5).Value = 118.99
$ws.Cells.Item(194, 6).Value = 1.3
$ws.Cells.Item(194, 7).Value = 37.1
$ws.Cells.Item(194, 8).Value = 126.71
$ws.Cells.Item(194, 9).Value = 6.35
$ws.Cells.Item(194, 10).Value = 95.01
$ws.Cells.Item(194, 11).Value = 92.66
$ws.Cells.Item(194, 12).Value = 1150.66
$ws.Cells.Item(194, 13).Value = 222.3
$ws.Cells.Item(194, 14).Value = 89.95
$ws.Cells.Item(194, 15).Value = 592.54
$ws.Cells.Item(194, 16).Value = 647.64
$ws.Cells.Item(194, 17).Value = 816.28
$ws.Cells.Item(194, 18).Value = 816.28
$ws.Cells.Item(194, 19).Value = 388.15
$ws.Cells.Item(194, 20).Value = 995.46
$ws.Cells.Item(194, 21).Value = 816.28
$ws.Cells.Item(194, 22).Value = 600.25
$ws.Cells.Item(194, 23).Value = 104.83
$ws.Cells.Item(194, 24).Value = 563.89
$ws.Cells.Item(194, 25).Value = 29.18
$ws.Cells.Item(194, 26).Value = 0.04
$ws.Cells.Item(194, 27).Value = 942.59
$ws.Cells.Item(194, 28).Value = 2.63
$ws.Cells.Item(194, 29).Value = 7.9
$ws.Cells.Item(194, 30).Value = 879.42
$ws.Cells.Item(194, 31).Value = 0.12
$ws.Cells.Item(194, 32).Value = 31
$ws.Cells.Item(194, 33).Value = 190.68
$ws.Cells.Item(194, 34).Value = 52.15
$ws.Cells.Item(194, 35).Value = 1107.27
$ws.Cells.Item(194, 36).Value = 91.95
$ws.Cells.Item(194, 37).Value = 197.43
$ws.Cells.Item(194, 38).Value = 8.25
$ws.Cells.Item(194, 39).Value = 0.22
$ws.Cells.Item(194, 40).Value = 34.01
$ws.Cells.Item(194, 41).Value = 14.54
$ws.Cells.Item(194, 42).Value = 16.05
$ws.Cells.Item(194, 43).Value = 39.54
$ws.Cells.Item(194, 44).Value = 19.03
$ws.Cells.Item(194, 45).Value = 105.64
$ws.Cells.Item(194, 46).Value = 54.23
$ws.Cells.Item(194, 47).Value = 147.84
$ws.Cells.Item(194, 48).Value = 0.02
$ws.Cells.Item(194, 49).Value = 217.67
$ws.Cells.Item(194, 50).Value = 195.17
$ws.Cells.Item(194, 51).Value = 11.27
$ws.Cells.Item(194, 52).Value = 10.89
$ws.Cells.Item(194, 53).Value = 0.06
$ws.Cells.Item(194, 54).Value = 4.78
$ws.Cells.Item(194, 55).Value = 252.18
$ws.Cells.Item(194, 56).Value = 1.92
$ws.Cells.Item(194, 57).Value = 0.68
$ws.Cells.Item(194, 58).Value = 7.33
$ws.Cells.Item(194, 59).Value = 126.41
$ws.Cells.Item(194, 60).Value = 206.86

# ---- Row 195 (08-10-2021) ----
Set-DateText 195 "08-10-2021"
$ws.Cells.Item(195, 2).Value = 24.09
$ws.Cells.Item(195, 3).Value = 813.62
$ws.Cells.Item(195, 4).Value = 195.81
$ws.Cells.Item(195, 5).Value = 118.6
$ws.Cells.Item(195, 6).Value = 1.3
$ws.Cells.Item(195, 7).Value = 36.98
$ws.Cells.Item(195, 8).Value = 126.44
$ws.Cells.Item(195, 9).Value = 6.32
$ws.Cells.Item(195, 10).Value = 95.03
$ws.Cells.Item(195, 11).Value = 92.71
$ws.Cells.Item(195, 12).Value = 1147.88
$ws.Cells.Item(195, 13).Value = 221.57
$ws.Cells.Item(195, 14).Value = 89.92
$ws.Cells.Item(195, 15).Value = 595.32
$ws.Cells.Item(195, 16).Value = 648.15
$ws.Cells.Item(195, 17).Value = 813.62
$ws.Cells.Item(195, 18).Value = 813.62
$ws.Cells.Item(195, 19).Value = 385.99
$ws.Cells.Item(195, 20).Value = 992.22
$ws.Cells.Item(195, 21).Value = 813.62
$ws.Cells.Item(195, 22).Value = 599.09
$ws.Cells.Item(195, 23).Value = 104.52
$ws.Cells.Item(195, 24).Value = 564.5
$ws.Cells.Item(195, 25).Value = 29.1
$ws.Cells.Item(195, 26).Value = 0.04
$ws.Cells.Item(195, 27).Value = 940.82
$ws.Cells.Item(195, 28).Value = 2.63
$ws.Cells.Item(195, 29).Value = 7.93
$ws.Cells.Item(195, 30).Value = 877.03
$ws.Cells.Item(195, 31).Value = 0.12
$ws.Cells.Item(195, 32).Value = 30.94
$ws.Cells.Item(195, 33).Value = 190.22
$ws.Cells.Item(195, 34).Value = 51.96
$ws.Cells.Item(195, 35).Value = 1108.32
$ws.Cells.Item(195, 36).Value = 91.71
$ws.Cells.Item(195, 37).Value = 198.86
$ws.Cells.Item(195, 38).Value = 8.22
$ws.Cells.Item(195, 39).Value = 0.22
$ws.Cells.Item(195, 40).Value = 33.9
$ws.Cells.Item(195, 41).Value = 14.47
$ws.Cells.Item(195, 42).Value = 16.09
$ws.Cells.Item(195, 43).Value = 39.46
$ws.Cells.Item(195, 44).Value = 18.84
$ws.Cells.Item(195, 45).Value = 105.59
$ws.Cells.Item(195, 46).Value = 54.52
$ws.Cells.Item(195, 47).Value = 147.89
$ws.Cells.Item(195, 48).Value = 0.02
$ws.Cells.Item(195, 49).Value = 216.97
$ws.Cells.Item(195, 50).Value = 194.55
$ws.Cells.Item(195, 51).Value = 11.34
$ws.Cells.Item(195, 52).Value = 10.88
$ws.Cells.Item(195, 53).Value = 0.06
$ws.Cells.Item(195, 54).Value = 4.78
$ws.Cells.Item(195, 55).Value = 252.38
$ws.Cells.Item(195, 56).Value = 1.91
$ws.Cells.Item(195, 57).Value = 0.68
$ws.Cells.Item(195, 58).Value = 7.29
$ws.Cells.Item(195, 59).Value = 126.08
$ws.Cells.Item(195, 60).Value = 205.64
